# Apply the "automatic update" of the artfynd export:
#  - Taxonsorteringsordning (col B) refreshed for a handful of rows
#  - Row 10 and Row 11 swap all their observation data (species, coords,
#    times, the "Ringhack" comment) because the underlying records were
#    re-sorted; only col B ends up with newly-recalculated sort keys
#    rather than simply following the swap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Addr,
        [string]$Text
    )
    # Force text so strings that look like dates/times/numbers aren't
    # auto-converted by Excel's smart entry, then strip the temporary
    # number-format override so the cell keeps the workbook's default
    # (unstyled) cell format, same as the source file.
    $r = $ws.Range($Addr)
    $r.NumberFormat = "@"
    $r.Value = $Text
    $r.ClearFormats()
}

# --- Row 3 ---------------------------------------------------------
$ws.Range("B3").Value = 79244

# --- Row 4 ---------------------------------------------------------
$ws.Range("B4").Value = 79244

# --- Row 10 ---------------------------------------------------------
# becomes what row 11 held, except column B gets a freshly recalculated
# sort key and the "Ringhack" public comment is dropped.
$ws.Range("A10").Value = 130894760
$ws.Range("B10").Value = 79244
$ws.Range("E10").Value = 6425
Set-TextValue "F10" "Garnlav"
Set-TextValue "G10" "Alectoria sarmentosa"
Set-TextValue "H10" "(Ach.) Ach."
$ws.Range("Q10").Value = 406786
$ws.Range("R10").Value = 7010890
Set-TextValue "Y10" "2026-01-18"
Set-TextValue "Z10" "14:31"
Set-TextValue "AA10" "2026-01-18"
Set-TextValue "AB10" "14:31"
$ws.Range("AC10").ClearContents()

# --- Row 11 ---------------------------------------------------------
# becomes what row 10 originally held, except column B gets a freshly
# recalculated sort key and the "Ringhack" public comment is gained.
$ws.Range("A11").Value = 130894767
$ws.Range("B11").Value = 57884
$ws.Range("E11").Value = 100109
Set-TextValue "F11" "Tretåig hackspett"
Set-TextValue "G11" "Picoides tridactylus"
Set-TextValue "H11" "(Linnaeus, 1758)"
$ws.Range("Q11").Value = 407194
$ws.Range("R11").Value = 7011100
Set-TextValue "Y11" "2026-01-21"
Set-TextValue "Z11" "12:26"
Set-TextValue "AA11" "2026-01-21"
Set-TextValue "AB11" "12:26"
Set-TextValue "AC11" "Ringhack"

# --- Row 14 ---------------------------------------------------------
$ws.Range("B14").Value = 79244

# --- Row 15 ---------------------------------------------------------
$ws.Range("B15").Value = 91829
